# "Cambridge University wording consistency"
# The workbook used "Cambridge University" in two places (Educational
# Background > Institution, and Organizations > Organization_name) while
# other cells referred to "University of Cambridge". Standardize both
# occurrences to "University of Cambridge".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Educational Background table: Institution column (row 19)
$ws.Range("B19").Value = "University of Cambridge"

# Organizations table: Organization_name column (row 27)
$ws.Range("A27").Value = "University of Cambridge"

# Leave the selection on B19, matching the saved worksheet view.
$ws.Range("B19").Select()
